$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A for the "#" numbering column, shifting existing
# Name/Issues columns one to the right.
$ws.Columns.Item(1).Insert()

# Header row
$ws.Cells.Item(1, 1).Value = "#"
$ws.Cells.Item(1, 2).Value = "Name"
$ws.Cells.Item(1, 3).Value = "Issues"

# Update existing Combank row text and add numbering
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "Combank"
$ws.Cells.Item(2, 3).Value = "Getting few null values for merchant details for general promos. Extracting it from promo title"

# New rows for DFCC, HNB, NDB - match the shared-string registration order
# of the original edit (bank names for rows 3 & 4, then the Issues text for
# row 3, then the final bank name for row 5).
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "DFCC"

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "HNB"

$ws.Cells.Item(3, 3).Value = "In promo details check how to make the main field (text before :) to bold like in the website"

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "NDB"

# Column widths to match target layout. Column B keeps the width it
# inherited from the original column A via the Insert() above, so it is
# left untouched here to avoid any re-quantization.
$ws.Columns.Item(1).ColumnWidth = 5.140625
$ws.Columns.Item(3).ColumnWidth = 68.28515625

# Selection as in target file
$ws.Range("B4").Select()
